$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update "last updated" timestamp string (A1) ---
$ws.Range("A1").Value = "Datos actualizados a 1 de Junio de 2020 a las 06:05"

# --- Swap country names between adjacent rows (data columns stay attached to the row) ---
# Belice / Santa Lucia swap (rows 201-202)
$ws.Range("A201").Value = "__TMP_SWAP_1__"
$ws.Range("A202").Value = "Belice"
$ws.Range("A201").Value = "Santa Lucia"

# Seychelles / Montserrat swap (rows 210-211)
$ws.Range("A210").Value = "__TMP_SWAP_2__"
$ws.Range("A211").Value = "Seychelles"
$ws.Range("A210").Value = "Montserrat"

# --- Update numeric data (India, row 10) ---
$ws.Cells.Item(10, 2).Value = 190622
$ws.Cells.Item(10, 3).Value = 13
$ws.Cells.Item(10, 4).Value = 91855
$ws.Cells.Item(10, 5).Value = 93359

# --- Update numeric data (Pakistan, row 21) ---
$ws.Cells.Item(21, 2).Value = 72460
$ws.Cells.Item(21, 3).Value = 2964
$ws.Cells.Item(21, 4).Value = 26083
$ws.Cells.Item(21, 5).Value = 44834
$ws.Cells.Item(21, 7).Value = 60
$ws.Cells.Item(21, 8).Value = 1543

# --- Update numeric data (Sri Lanka, row 101) ---
$ws.Cells.Item(101, 5).Value = 821
$ws.Cells.Item(101, 7).Value = 1
$ws.Cells.Item(101, 8).Value = 11

# --- Update numeric data (Mongolia, row 161) ---
$ws.Cells.Item(161, 2).Value = 185
$ws.Cells.Item(161, 3).Value = 6
$ws.Cells.Item(161, 5).Value = 141

# --- Update numeric data for swapped rows 201/202 (Santa Lucia / Belice) ---
$ws.Cells.Item(201, 4).Value = 18
$ws.Cells.Item(201, 8).Value = 0
$ws.Cells.Item(202, 4).Value = 16
$ws.Cells.Item(202, 8).Value = 2

# --- Update numeric data for swapped rows 210/211 (Montserrat / Seychelles) ---
$ws.Cells.Item(210, 4).Value = 10
$ws.Cells.Item(210, 8).Value = 1
$ws.Cells.Item(211, 4).Value = 11
$ws.Cells.Item(211, 8).Value = 0
